$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.404.72'
$ws.Range("E2").Value = '  -2.87%  '

$ws.Range("D3").Value = '3.156.49'
$ws.Range("E3").Value = '  -4.14%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '527.87'
$ws.Range("E5").Value = '  -5.11%  '

$ws.Range("D6").Value = '134.51'
$ws.Range("E6").Value = '  -4.83%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.153.58'
$ws.Range("E8").Value = '  -4.21%  '

$ws.Range("D9").Value = '0.445'
$ws.Range("E9").Value = '  -4.57%  '

$ws.Range("E10").Value = '  -7.54%  '

$ws.Range("E11").Value = '  -7.70%  '

$ws.Range("E12").Value = '  -7.17%  '

$ws.Range("D13").Value = '3.695.63'
$ws.Range("E13").Value = '  -4.08%  '

$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("E15").Value = '  -5.06%  '

$ws.Range("D16").Value = '3.156.29'
$ws.Range("E16").Value = '  -3.90%  '

$ws.Range("D17").Value = '58.359.33'
$ws.Range("E17").Value = '  -3.01%  '

$ws.Range("E18").Value = '  -6.97%  '

$ws.Range("E19").Value = '  -5.09%  '

$ws.Range("D20").Value = '13.07'
$ws.Range("E20").Value = '  -5.10%  '

$ws.Range("E21").Value = '  -6.72%  '

$ws.Range("D22").Value = '343.16'
$ws.Range("E22").Value = '  -7.79%  '

$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("E24").Value = '  -4.07%  '

$ws.Range("D25").Value = '67.14'
$ws.Range("E25").Value = '  -7.15%  '

$ws.Range("D26").Value = '3.284.61'
$ws.Range("E26").Value = '  -3.85%  '

$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").Value = '0.0₃0953'
$ws.Range("E28").Value = '  -7.08%  '

$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").Value = '6.84'
$ws.Range("E30").Value = '  -3.09%  '

$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.88'
$ws.Range("E32").Value = '  -6.75%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.25'
$ws.Range("E33").Value = '  +2.01%  '

$ws.Range("E34").Value = '  -6.90%  '

$ws.Range("D35").Value = '21.50'
$ws.Range("E35").Value = '  -4.63%  '

$ws.Range("D36").Value = '4.85'
$ws.Range("E36").Value = '  -4.18%  '

$ws.Range("D37").Value = '159.54'
$ws.Range("E37").Value = '  -3.86%  '

$ws.Range("E38").Value = '  -5.42%  '

$ws.Range("E39").Value = '  -9.32%  '

$ws.Range("E40").Value = '  -4.67%  '

$ws.Range("D41").Value = '3.185.89'
$ws.Range("E41").Value = '  -4.07%  '

$ws.Range("D42").Value = '40.45'
$ws.Range("E42").Value = '  -2.68%  '

$ws.Range("D43").Value = '24.02'
$ws.Range("E43").Value = '  -6.99%  '

$ws.Range("E45").Value = '  -2.62%  '

$ws.Range("D46").Value = '3.94'
$ws.Range("E46").Value = '  -3.83%  '

$ws.Range("D47").Value = '0.999'

$ws.Range("E48").Value = '  -6.92%  '

$ws.Range("D49").Value = '2.291.02'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").Value = '6.18'
$ws.Range("E50").Value = '  -2.56%  '

$ws.Range("D51").Value = '20.72'
$ws.Range("E51").Value = '  -3.67%  '
